# Update cryptos list with refreshed market data (prices / 1h volume change).
# Rows 14 and 15 also swap coin identity: WrappedliquidstakedEther2.0 moves
# above ShibaInu in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.376.66"
$ws.Range("E2").Value = "  -1.32%  "

# Row 3
$ws.Range("D3").Value = "2.426.10"
$ws.Range("E3").Value = "  -2.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.22%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("E8").Value = "  -0.23%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.162"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.33%  "

# Row 10
$ws.Range("E10").Value = "  -1.80%  "

# Row 11
$ws.Range("E11").Value = "  -0.93%  "

# Row 12
$ws.Range("E12").Value = "  -5.77%  "

# Row 13
$ws.Range("D13").Value = "68.240.13"
$ws.Range("E13").Value = "  -1.33%  "

# Rows 14-15: WrappedliquidstakedEther2.0 and ShibaInu swap ranking positions
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.871.93"
$ws.Range("E14").Value = "  -1.16%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000173"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.47%  "

# Row 17
$ws.Range("D17").Value = "2.426.08"
$ws.Range("E17").Value = "  -1.30%  "

# Row 18
$ws.Range("E18").Value = "  -3.69%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "334.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.75%  "

# Row 20
$ws.Range("E20").Value = "  -2.42%  "

# Row 21
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.46%  "

# Row 23
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "

# Row 26
$ws.Range("D26").Value = "2.547.36"
$ws.Range("E26").Value = "  -2.20%  "

# Row 27
$ws.Range("E27").Value = "  -1.08%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0812"
$ws.Range("E28").Value = "  -1.89%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "425.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.11%  "

# Row 32
$ws.Range("E32").Value = "  -0.88%  "

# Row 33
$ws.Range("E33").Value = "  -1.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.46%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.23%  "

# Row 36
$ws.Range("E36").Value = "  -0.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.104"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.14%  "

# Row 39
$ws.Range("E39").Value = "  -2.06%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.50%  "

# Row 41
$ws.Range("E41").Value = "  -0.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.59%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "132.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.60%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.60%  "

# Row 45
$ws.Range("E45").Value = "  -1.14%  "

# Row 46
$ws.Range("E46").Value = "  -0.81%  "

# Row 47
$ws.Range("E47").Value = "  -1.68%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.556"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.63%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0914"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.48%  "

# Row 50
$ws.Range("E50").Value = "  -0.10%  "

# Row 51
$ws.Range("E51").Value = "  -2.61%  "

